# The sheet previously carried a (stripped-out) ~356k row dataset; the
# fixture was trimmed down to a single representative data row so the
# file no longer needs to live in git-lfs. Re-create that single row:
# year 2005 in A2, leaving the cursor on the cell just typed into.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 2005
[void]$ws.Range("A2").Select()
